$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

$ws.Range("H52").Value = 1618
$ws.Range("I52").Value = 1618
$ws.Range("K52").Value = 4854
$ws.Range("M52").Value = -4694

$ws.Range("H58").Value = 3134.3333
$ws.Range("J58").Value = 8888
$ws.Range("L58").Value = 26664
$ws.Range("N58").Value = -26964

$ws.Range("H92").Value = 77778580
$ws.Range("I92").Value = 13889886
$ws.Range("J92").Value = 333333340
$ws.Range("K92").Value = 13889886
$ws.Range("L92").Value = 333333340
$ws.Range("M92").Value = -13888638
$ws.Range("N92").Value = -333335836

$ws.Range("H101").Value = 890
$ws.Range("I101").Value = 900
$ws.Range("J101").Value = 883.3333
$ws.Range("K101").Value = 2700
$ws.Range("L101").Value = 2649.9999
$ws.Range("M101").Value = -1078
$ws.Range("N101").Value = -5893.9999

$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws.Range("H107").Value = 62500976
$ws.Range("I107").Value = 62500976
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 62500976
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -62499056
$ws.Range("N107").ClearContents()

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws.Range("H112").Value = 21979304
$ws.Range("J112").Value = 24846058
$ws.Range("L112").Value = 74538174
$ws.Range("N112").Value = -74540390

$ws.Range("H115").Value = 340
$ws.Range("I115").Value = 340
$ws.Range("K115").Value = 1020
$ws.Range("M115").Value = 547

$ws.Range("H136").Value = 49106.25
$ws.Range("J136").Value = 49106.25
$ws.Range("L136").Value = 49106.25
$ws.Range("N136").Value = -59306.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 3290.4
$ws.Range("I33").Value = 1863
$ws.Range("J33").Value = 9000
$ws.Range("K33").Value = 1863
$ws.Range("L33").Value = 9000
$ws.Range("M33").Value = -1534
$ws.Range("N33").Value = -9658

$ws.Range("H133").Value = 32875
$ws.Range("J133").Value = 32875
$ws.Range("L133").Value = 32875
$ws.Range("N133").Value = -37935

$ws.Range("H134").Value = 50293.332
$ws.Range("J134").Value = 50293.332
$ws.Range("L134").Value = 50293.332
$ws.Range("N134").Value = -60433.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 295.5
$ws.Range("I22").Value = 49.2
$ws.Range("J22").Value = 471.42856
$ws.Range("K22").Value = 49.2
$ws.Range("L22").Value = 471.42856
$ws.Range("M22").Value = 123.8
$ws.Range("N22").Value = -817.4285600000001

$ws.Range("H132").Value = 51372
$ws.Range("J132").Value = 51372
$ws.Range("L132").Value = 51372
$ws.Range("N132").Value = -61492

$ws.Range("H133").Value = 43696.668
$ws.Range("J133").Value = 44436
$ws.Range("L133").Value = 44436
$ws.Range("N133").Value = -54556

$ws.Range("H138").Value = 60560
$ws.Range("J138").Value = 60560
$ws.Range("L138").Value = 60560
$ws.Range("N138").Value = -70840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 41819
$ws.Range("J3").Value = 100003
$ws.Range("L3").Value = 100003
$ws.Range("N3").Value = -100229

$ws.Range("H31").Value = 392679.53
$ws.Range("I31").Value = 1887.3846
$ws.Range("J31").Value = 580838.7
$ws.Range("K31").Value = 1887.3846
$ws.Range("L31").Value = 580838.7
$ws.Range("M31").Value = -1592.3846
$ws.Range("N31").Value = -581428.7

$ws.Range("H32").Value = 2003.3334
$ws.Range("I32").Value = 2003.3334
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2003.3334
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1687.3334
$ws.Range("N32").ClearContents()

$ws.Range("H34").Value = 392679.53
$ws.Range("I34").Value = 1887.3846
$ws.Range("J34").Value = 580838.7
$ws.Range("K34").Value = 1887.3846
$ws.Range("L34").Value = 580838.7
$ws.Range("M34").Value = -1685.3846
$ws.Range("N34").Value = -581242.7

$ws.Range("H99").Value = 10419277
$ws.Range("I99").Value = 2022.4
$ws.Range("J99").Value = 17860174
$ws.Range("K99").Value = 2022.4
$ws.Range("L99").Value = 17860174
$ws.Range("M99").Value = -524.4000000000001
$ws.Range("N99").Value = -17863170

$ws.Range("H126").Value = 10419277
$ws.Range("I126").Value = 2022.4
$ws.Range("J126").Value = 17860174
$ws.Range("K126").Value = 6067.200000000001
$ws.Range("L126").Value = 53580522
$ws.Range("M126").Value = -3597.200000000001
$ws.Range("N126").Value = -53585462

$ws.Range("H133").Value = 51575.6
$ws.Range("J133").Value = 51575.6
$ws.Range("L133").Value = 51575.6
$ws.Range("N133").Value = -56635.6

$ws.Range("H138").Value = 58397.5
$ws.Range("J138").Value = 58397.5
$ws.Range("L138").Value = 58397.5
$ws.Range("N138").Value = -68677.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1159.3334
$ws.Range("I2").Value = 2025.8
$ws.Range("J2").Value = 76.25
$ws.Range("K2").Value = 12154.8
$ws.Range("L2").Value = 457.5
$ws.Range("M2").Value = -12041.8
$ws.Range("N2").Value = -683.5

$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()

$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()

$ws.Range("H37").Value = 41000
$ws.Range("J37").Value = 41000
$ws.Range("L37").Value = 123000
$ws.Range("N37").Value = -123224

$ws.Range("H113").Value = 3750490.2
$ws.Range("I113").Value = 4545909
$ws.Range("K113").Value = 13637727
$ws.Range("M113").Value = -13635557

$ws.Range("H131").Value = 2703854.8
$ws.Range("I131").Value = 5882748.5
$ws.Range("J131").Value = 1795
$ws.Range("K131").Value = 17648245.5
$ws.Range("L131").Value = 5385
$ws.Range("M131").Value = -17643205.5
$ws.Range("N131").Value = -15465

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1331.6666
$ws.Range("I97").Value = 1369.6428
$ws.Range("J97").Value = 800
$ws.Range("K97").Value = 1369.6428
$ws.Range("L97").Value = 800
$ws.Range("M97").Value = -873.6428000000001
$ws.Range("N97").Value = -1792

$ws.Range("H102").Value = 3410.7778
$ws.Range("I102").Value = 2874.75
$ws.Range("J102").Value = 3839.6
$ws.Range("K102").Value = 2874.75
$ws.Range("L102").Value = 3839.6
$ws.Range("M102").Value = -1252.75
$ws.Range("N102").Value = -7083.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 1780.25
$ws.Range("I35").Value = 1780.25
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1780.25
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -1444.25
$ws.Range("N35").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 67169.836
$ws.Range("I3").Value = 1003
$ws.Range("J3").Value = 80403.2
$ws.Range("K3").Value = 1003
$ws.Range("L3").Value = 80403.2
$ws.Range("M3").Value = -889
$ws.Range("N3").Value = -80631.2

$ws.Range("H32").Value = 12333.333
$ws.Range("J32").Value = 11000
$ws.Range("L32").Value = 11000
$ws.Range("N32").Value = -11634

$ws.Range("H135").Value = 39483.855
$ws.Range("J135").Value = 39483.855
$ws.Range("L135").Value = 39483.855
$ws.Range("N135").Value = -49623.855

$ws.Range("H141").Value = 71487.336
$ws.Range("J141").Value = 71487.336
$ws.Range("L141").Value = 71487.336
$ws.Range("N141").Value = -81847.336
